$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row text (A1:C1 keep their content, shifted/renamed)
$ws.Range("A1").Value = "Descricao do Insumo"
$ws.Range("B1").Value = "Valor unitario"
$ws.Range("C1").Value = "Problema(s)"

# Remove the now-unused trailing columns D:G from the header row
$ws.Range("D1:G1").ClearContents()
